# Update the RAD test result timestamps in column B (rows 2-4) on Sheet1.
# These cells log the "Date" the Katalon RAD test case executed; the commit
# re-runs the test suite, producing fresh timestamps for MRF test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Oct 02 16:45:44 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 16:45:57 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 16:46:10 EDT 2023"
